# Generate Report for Handback
#
# - Status text moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   (shared by the Overview tab and both language tabs).
# - Each language tab grows two new columns of data for row 2/3:
#     F = Latest Target File   (same display text/link as the source file, col A)
#     G = Latest Handback File (same display text/link style as col D)
# - "Latest Handback DateTime" (col H) is stamped with the real handback time
#   (zh-cn finishes first, de-de a few seconds later).

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # OLE (BGR) form of RGB(0x64,0x95,0xED) used by the workbook's HyperLink style

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2      # xlUnderlineStyleSingle
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn tab
# ---------------------------------------------------------------------------

$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("F2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $zhcn.Range("F2")

$zhcn.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fd0c922429f15ebf7c308b546ab7ed307f2c645/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
Style-AsHyperlink $zhcn.Range("G2")

$zhcn.Range("F3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $zhcn.Range("F3")

$zhcn.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fd0c922429f15ebf7c308b546ab7ed307f2c645/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
Style-AsHyperlink $zhcn.Range("G3")

$zhcn.Range("H2").Value = "2016-03-24 20:41:20"
$zhcn.Range("H3").Value = "2016-03-24 20:41:20"

# ---------------------------------------------------------------------------
# 3. de-de tab
# ---------------------------------------------------------------------------

$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("F2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $dede.Range("F2")

$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f983606bf0385f1892ead680a99baafc3ba145b4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
Style-AsHyperlink $dede.Range("G2")

$dede.Range("F3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/a.md", [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $dede.Range("F3")

$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f983606bf0385f1892ead680a99baafc3ba145b4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
Style-AsHyperlink $dede.Range("G3")

$dede.Range("H2").Value = "2016-03-24 20:41:29"
$dede.Range("H3").Value = "2016-03-24 20:41:29"
